# Fix analytics API errors and remove duplicate offer details

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix application status for row 2 (John Doe) ---
$ws.Range("N2").Value = "Rejected"

# --- Fill in missing Recommendation field in Round 1 Remarks JSON for row 2 ---
$ws.Range("AA2").Value = '{"Communication":"good ","Technical Assessment":"good","Problem-Solving":"","Overall Potential":"","Recommendation":"Proceed Round 2"}'

# --- Append new candidate rows (18 and 19) ---
# Columns that hold numeric-looking text (phone numbers / CTC figures) are
# pre-formatted as Text so Excel keeps them as strings (preserving leading
# zeros etc.) instead of silently coercing them to numbers.

# Row 18: JAGADEESH M
$ws.Cells.Item(18, 1).Value = ""
$ws.Cells.Item(18, 2).Value = "JAGADEESH M"
$ws.Cells.Item(18, 3).Value = "jagadeesh19ct11@gmail.com"
$ws.Range("D18").NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "06381813711"
$ws.Cells.Item(18, 5).Value = "Site Reliability Engineer"
$ws.Cells.Item(18, 6).Value = "Data Engineer"
$ws.Cells.Item(18, 7).Value = "nill"
$ws.Cells.Item(18, 8).Value = "Coimbatore"
$ws.Range("I18").NumberFormat = "@"
$ws.Cells.Item(18, 9).Value = "0"
$ws.Range("J18").NumberFormat = "@"
$ws.Cells.Item(18, 10).Value = "400000"
$ws.Cells.Item(18, 11).Value = "0-1 years"
$ws.Cells.Item(18, 12).Value = "Immediate"
$ws.Cells.Item(18, 13).Value = ""
$ws.Cells.Item(18, 14).Value = "Rejected"
$ws.Cells.Item(18, 15).Value = ""
$ws.Cells.Item(18, 16).Value = ""
$ws.Cells.Item(18, 17).Value = "No"
$ws.Cells.Item(18, 18).Value = "Yes"
$ws.Cells.Item(18, 19).Value = "No"
$ws.Range("T18").NumberFormat = "@"
$ws.Cells.Item(18, 20).Value = "0"
$ws.Cells.Item(18, 21).Value = "Coimbatore"
$ws.Cells.Item(18, 22).Value = "nill"
$ws.Cells.Item(18, 23).Value = ""
$ws.Cells.Item(18, 24).Value = ""
$ws.Cells.Item(18, 25).Value = ""
$ws.Cells.Item(18, 26).Value = ""
$ws.Cells.Item(18, 27).Value = ""
$ws.Cells.Item(18, 28).Value = ""
$ws.Cells.Item(18, 29).Value = ""
$ws.Cells.Item(18, 30).Value = ""
$ws.Cells.Item(18, 31).Value = "Yes"
$ws.Cells.Item(18, 32).Value = ""

# Row 19: ashwin
$ws.Cells.Item(19, 1).Value = ""
$ws.Cells.Item(19, 2).Value = "ashwin"
$ws.Cells.Item(19, 3).Value = "ashlog559@gmail.com"
$ws.Range("D19").NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "0123456789"
$ws.Cells.Item(19, 5).Value = "Senior Site Reliability Engineer"
$ws.Cells.Item(19, 6).Value = "Senior Software Engineer"
$ws.Cells.Item(19, 7).Value = "nill"
$ws.Cells.Item(19, 8).Value = "Chennai"
$ws.Range("I19").NumberFormat = "@"
$ws.Cells.Item(19, 9).Value = "100000000"
$ws.Range("J19").NumberFormat = "@"
$ws.Cells.Item(19, 10).Value = "999999999999998"
$ws.Cells.Item(19, 11).Value = "2-3 years"
$ws.Cells.Item(19, 12).Value = "Immediate"
$ws.Cells.Item(19, 13).Value = ""
$ws.Cells.Item(19, 14).Value = "Rejected"
$ws.Cells.Item(19, 15).Value = ""
$ws.Cells.Item(19, 16).Value = ""
$ws.Cells.Item(19, 17).Value = "Yes"
$ws.Cells.Item(19, 18).Value = "No"
$ws.Cells.Item(19, 19).Value = "Yes"
$ws.Range("T19").NumberFormat = "@"
$ws.Cells.Item(19, 20).Value = "999999999999"
$ws.Cells.Item(19, 21).Value = "Others"
$ws.Cells.Item(19, 22).Value = "nill"
$ws.Cells.Item(19, 23).Value = ""
$ws.Cells.Item(19, 24).Value = ""
$ws.Cells.Item(19, 25).Value = ""
$ws.Cells.Item(19, 26).Value = ""
$ws.Cells.Item(19, 27).Value = ""
$ws.Cells.Item(19, 28).Value = ""
$ws.Cells.Item(19, 29).Value = ""
$ws.Cells.Item(19, 30).Value = ""
$ws.Cells.Item(19, 31).Value = "Yes"
$ws.Cells.Item(19, 32).Value = ""
